# Update workbook/worksheet for SCD0012-002
# (renumbered from SCD0212 -> SCD0012, TC_ID updated from DGS-227 -> SCD0012-002,
#  column B widened to fit the new TC_ID text, and the view/selection reset)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet tab: SCD0212 -> SCD0012
$ws.Name = "SCD0012"

# 2. Update the TC_ID cell (B2) from "DGS-227" to "SCD0012-002"
$ws.Range("B2").Value = "SCD0012-002"

# 3. Widen column B so the longer TC_ID value fits (target stored width 12.42578125;
#    the host rounds ColumnWidth to the nearest 1/6 character, so feed it the input
#    that rounds to the closest achievable stored width)
$ws.Columns("B").ColumnWidth = 11.666666666666666

# 4. Reset the view: scroll back to show column A, and select B3 instead of L3
$ws.Range("B3").Select()
